$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row and the "SC 92" row entirely (rows shift up).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Update/clear the E-column (column D label) imputed values.
$ws.Range("E3").Value = -5.7

$ws.Range("E5").ClearContents()
$ws.Range("E5").Style = "Normal"

$ws.Range("E21").Value = -8.699999999999999

$ws.Range("E23").ClearContents()
$ws.Range("E23").Style = "Normal"

$ws.Range("E32").Value = -6.4
